# Generate Report for Handback
#
# Marks the two tracked files as handed back (in sync with en-US), records
# the handback timestamps, and fills in the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns on the
# zh-cn and de-de report sheets (previously blank placeholders).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "4e81a5eb-e5c7-4a2d-89f8-7817de87db74.md"
$mdFile2 = "bd5c5283-3373-40e9-94c5-2bf992a88840.md"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a253a538104a8e5c2dd1cc991d2973c4fc5fa9b/e2e/4e81a5eb-e5c7-4a2d-89f8-7817de87db74.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a253a538104a8e5c2dd1cc991d2973c4fc5fa9b/e2e/bd5c5283-3373-40e9-94c5-2bf992a88840.md"

$zhTargetFile1 = "4e81a5eb-e5c7-4a2d-89f8-7817de87db74.92b1504521bb9b073ba69a524a75df6cc6e9e78b.zh-cn.xlf"
$zhTargetFile2 = "bd5c5283-3373-40e9-94c5-2bf992a88840.e9a098edf0ae5b8e952b76119eca0fec99914762.zh-cn.xlf"
$deTargetFile1 = "4e81a5eb-e5c7-4a2d-89f8-7817de87db74.92b1504521bb9b073ba69a524a75df6cc6e9e78b.de-de.xlf"
$deTargetFile2 = "bd5c5283-3373-40e9-94c5-2bf992a88840.e9a098edf0ae5b8e952b76119eca0fec99914762.de-de.xlf"

$zhHandbackTime = "2016-08-28 21:00:12"
$deHandbackTime = "2016-08-28 21:00:20"

# Column widths: the quantised ColumnWidth setter on this runtime snaps to
# 1/6-character steps, so these are the closest reachable values to the
# wider columns in the handed-back report (~29.98 / 40 stored width).
$wideStatusWidth = 175 / 6    # -> stored width 30 (closest reachable to 29.9777)
$wideLinkWidth    = 235 / 6   # -> stored width 40 (exact)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status text + widened zh-cn/de-de columns ---------
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth

# --- zh-cn sheet ---------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url1, "", "", $mdFile1)
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Range("J2").Value = $zhTargetFile1
$wsZhCn.Range("K2").Value = $zhHandbackTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url2, "", "", $mdFile2)
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Range("J3").Value = $zhTargetFile2
$wsZhCn.Range("K3").Value = $zhHandbackTime

$wsZhCn.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $wideLinkWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideLinkWidth

# --- de-de sheet -----------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url1, "", "", $mdFile1)
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Range("J2").Value = $deTargetFile1
$wsDeDe.Range("K2").Value = $deHandbackTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url2, "", "", $mdFile2)
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Range("J3").Value = $deTargetFile2
$wsDeDe.Range("K3").Value = $deHandbackTime

$wsDeDe.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $wideLinkWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideLinkWidth
